$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.724.93'
$ws.Range("E2").Value = '  +0.42%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.699.91'
$ws.Range("E3").Value = '  +0.21%  '

# Row 4
$ws.Range("E4").Value = '  +0.34%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.16'
$ws.Range("E5").Value = '  -0.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.39%  '

# Row 7
$ws.Range("E7").Value = '  -0.50%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4036'
$ws.Range("E8").Value = '  +0.41%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.526'
$ws.Range("E9").Value = '  -0.58%  '

# Row 10
$ws.Range("B10").Value = 'BinanceUSD'
$ws.Range("C10").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.005'
$ws.Range("E10").Value = '  +0.31%  '

# Row 11
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.58'
$ws.Range("E11").Value = '  +0.43%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08859'
$ws.Range("E12").Value = '  +0.92%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.404'
$ws.Range("E13").Value = '  +2.13%  '

# Row 14
$ws.Range("E14").Value = '  +1.64%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.176'
$ws.Range("E15").Value = '  +7.46%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001323'
$ws.Range("E16").Value = '  +0.38%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.712.71'
$ws.Range("E17").Value = '  +1.38%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.70'
$ws.Range("E18").Value = '  -1.33%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07049'
$ws.Range("E19").Value = '  +0.53%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.69'
$ws.Range("E20").Value = '  -0.08%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.086'
$ws.Range("E21").Value = '  +2.76%  '

# Row 22
$ws.Range("E22").Value = '  +0.65%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.78'
$ws.Range("E23").Value = '  +5.06%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.713.58'
$ws.Range("E24").Value = '  +0.41%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.135'
$ws.Range("E25").Value = '  +2.17%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.370'
$ws.Range("E26").Value = '  +1.42%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.70'
$ws.Range("E27").Value = '  +1.56%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '163.14'
$ws.Range("E28").Value = '  +2.03%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.719'
$ws.Range("E29").Value = '  +16.34%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.80'
$ws.Range("E30").Value = '  +1.32%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.183'
$ws.Range("E31").Value = '  -0.11%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09039'
$ws.Range("E32").Value = '  +6.28%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.645'
$ws.Range("E33").Value = '  +3.27%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.065'
$ws.Range("E34").Value = '  -3.16%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.983'
$ws.Range("E35").Value = '  +0.84%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.12'
$ws.Range("E36").Value = '  -3.55%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2755'
$ws.Range("E37").Value = '  +0.52%  '

# Row 38
$ws.Range("E38").Value = '  -0.31%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02772'
$ws.Range("E39").Value = '  +0.67%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09146'
$ws.Range("E40").Value = '  +1.40%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.464'
$ws.Range("E41").Value = '  -0.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7675'
$ws.Range("E42").Value = '  -0.46%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.87'
$ws.Range("E43").Value = '  +2.41%  '

# Row 44
$ws.Range("E44").Value = '  -0.47%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.571'
$ws.Range("E45").Value = '  +1.61%  '

# Row 46
$ws.Range("E46").Value = '  -0.18%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("E47").Value = '  +0.49%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.335'
$ws.Range("E48").Value = '  -0.47%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.08'
$ws.Range("E49").Value = '  -0.54%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '90.87'
$ws.Range("E50").Value = '  +3.03%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07989'
$ws.Range("E51").Value = '  -0.60%  '
